$d = $word.ActiveDocument

# Locate the target paragraph: "ngShow, ngHide (tags pequenas) e ngIf ..."
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "ngShow, ngHide*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'ngShow, ngHide' paragraph"
}

$rng = $target.Range

# Replace the single paragraph with three paragraphs:
#  1) the original ngShow/ngHide/ngIf paragraph, with the paragraph-mark
#     bold removed, "(tags pequenas)" merged into one run, and the
#     "(melhor para performance - associar ngsource)" merged into one run.
#  2) a new "ngInclude: ..." paragraph, which now also owns the _GoBack
#     bookmark.
#  3) a new "ngRequired ..." paragraph (bold paragraph mark restored).
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr><w:jc w:val="both"/></w:pPr>
  <w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">ngShow, ngHide </w:t></w:r>
  <w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">(tags pequenas) </w:t></w:r>
  <w:r><w:rPr><w:b/></w:rPr><w:t>e ngIf</w:t></w:r>
  <w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t>(melhor para performance - associar ngsource)</w:t></w:r>
  <w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r>
  <w:r><w:t xml:space="preserve"> Exibindo um elemento condicionalmente.</w:t></w:r>
  <w:r><w:t xml:space="preserve"> ngIf interage com a DOM já os outros não.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:jc w:val="both"/></w:pPr>
  <w:r><w:t>ngInclude: Incluir conteúdo dinamicamente.</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p>
  <w:pPr><w:jc w:val="both"/><w:rPr><w:b/></w:rPr></w:pPr>
  <w:r><w:rPr><w:b/></w:rPr><w:t>ngRequired</w:t></w:r>
  <w:r><w:t xml:space="preserve"> Define um determinado campo como obrigatório.</w:t></w:r>
  <w:r><w:t xml:space="preserve"> Com </w:t></w:r>
  <w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">$invalid </w:t></w:r>
  <w:r><w:t>para desabilitar &#8220;botão&#8221;</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rng.InsertXML($xml)
